# Mengubah index pasien menjadi id pasien
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace numeric patient index (1-4) with patient id strings (P1-P4) in A2:A5
$ws.Range("A2").Value = "P1"
$ws.Range("A3").Value = "P2"
$ws.Range("A4").Value = "P3"
$ws.Range("A5").Value = "P4"

# Update selection to match the new active selection of A2:A5
$ws.Range("A2:A5").Select()
